# Apply updated dSF (column F) values per the "repull data, push all data, mean calculation" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new value for column F (dSF)
$updates = @{
    2  = -8
    3  = -6
    4  = -2
    5  = -5
    7  = 1
    8  = -3
    9  = 2
    10 = -1
    11 = -1
    13 = -6
    17 = 1
    18 = -4
    19 = 3
    23 = -8
    24 = -2
    25 = -1
    27 = -6
    30 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
